$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# Bring in the cell formatting (styles) for the new rows 69-72 by
# copying the formats from existing rows that already use the same
# style combination used by the new "Organization" rule rows.
# ------------------------------------------------------------------
$ws.Range("B20:G20").Copy()
$ws.Range("B69:G69").PasteSpecial(-4122)

$ws.Range("B32:G32").Copy()
$ws.Range("B70:G70").PasteSpecial(-4122)

$ws.Range("B32:G32").Copy()
$ws.Range("B71:G71").PasteSpecial(-4122)

$ws.Range("B65:G65").Copy()
$ws.Range("B72:G72").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Rows.Item(69).RowHeight = 30
$ws.Rows.Item(70).RowHeight = 30
$ws.Rows.Item(71).RowHeight = 30
$ws.Rows.Item(72).RowHeight = 30

# ------------------------------------------------------------------
# Populate the new "Organization" access-control-list rows.
# ------------------------------------------------------------------
$ws.Range("C69").Value = "ORGANIZATION"
$ws.Range("B70").Value = "Organization – default read access"
$ws.Range("B71").Value = "Organization – Only participants can save"
$ws.Range("B69").Value = "Organization – Lockout No Access Users"
$ws.Range("G71").Value = "grant save to owner, owning group, collaborator"
$ws.Range("G70").Value = "grant read to owner, owning group, collaborator, reader, *"
$ws.Range("B72").Value = "Organization – Restricted Flag"

$ws.Range("C70").Value = "ORGANIZATION"
$ws.Range("C71").Value = "ORGANIZATION"
$ws.Range("C72").Value = "ORGANIZATION"
$ws.Range("D72").Value = "restricted"
$ws.Range("G69").Value = "mandatory deny read to No Access"
$ws.Range("G72").Value = "deny read to *"

# ------------------------------------------------------------------
# Update the view so the new rows are visible / selected, matching
# the author's final cursor position.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 63
$ws.Range("G69").Select()
